$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update cell values on the "Priority" (C) column and add a missing "Comments" (D) entry ---
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = "DONE"

$ws.Range("C4").Value = 4

$ws.Range("C5").Value = 10

$ws.Range("C11").Value = 10

# --- Clear the existing AutoFilter criteria (was filtering Priority = 1) and unhide all rows ---
# Remove the criteria from column C (3rd column of the filtered range) while keeping the filter on.
$ws.Range("A2:D15").AutoFilter(3)

# Turn the AutoFilter off completely, then re-apply it across the full (now 16-row) data range
# with no active criteria, matching the target ref="A2:D16" with no <filterColumn>.
$ws.AutoFilterMode = $false
$ws.Range("A2:D16").AutoFilter()

# --- Update the active selection to C4 ---
$ws.Range("C4").Select()

# --- Keep the hidden "_FilterDatabase" defined name in sync with the new filter range ---
$fdb = $wb.Names.Item("Sheet1!_FilterDatabase")
$fdb.RefersTo = "=Sheet1!`$A`$2:`$D`$16"
